$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (sCs / Gdnf / Ret / ECs) ---
$ws.Cells.Item(2,5).Value  = 3                     # E2
$ws.Cells.Item(2,7).Value  = 1.995314333333334      # G2
$ws.Cells.Item(2,8).Value  = 5.985943000000001      # H2
$ws.Cells.Item(2,11).Value = 3                      # K2
$ws.Cells.Item(2,13).Value = 3.632885               # M2
$ws.Cells.Item(2,14).Value = 10.898655              # N2
$ws.Cells.Item(2,15).Value = 0.4372849566404539     # O2
$ws.Cells.Item(2,16).Value = 0.4372849566404539     # P2
$ws.Cells.Item(2,17).Value = 7.248747511851668      # Q2
$ws.Cells.Item(2,18).Value = 65.23872760666501      # R2
$ws.Cells.Item(2,19).Value = 0.4372849566404539     # S2
$ws.Cells.Item(2,20).Value = 0.4372849566404539     # T2

# --- Update row 3 (sCs / Gdnf / Ret / FAPs) ---
$ws.Cells.Item(3,5).Value  = 3                      # E3
$ws.Cells.Item(3,7).Value  = 1.995314333333334      # G3
$ws.Cells.Item(3,8).Value  = 5.985943000000001      # H3
$ws.Cells.Item(3,11).Value = 3                      # K3
$ws.Cells.Item(3,13).Value = 3.884996               # M3
$ws.Cells.Item(3,14).Value = 11.654988              # N3
$ws.Cells.Item(3,15).Value = 0.4676311822169809     # O3
$ws.Cells.Item(3,16).Value = 0.4676311822169809     # P3
$ws.Cells.Item(3,17).Value = 7.751788203742668      # Q3
$ws.Cells.Item(3,18).Value = 69.76609383368401      # R3
$ws.Cells.Item(3,19).Value = 0.4676311822169809     # S3
$ws.Cells.Item(3,20).Value = 0.4676311822169809     # T3

# --- Row 4 used to be (sCs / Gdnf / Ret / sCs); it now becomes (sCs / Gdnf / Ret / M1) ---
# with a brand-new set of values, and the old row 4 data moves down to become the new row 5.
$ws.Cells.Item(4,4).Value  = "M1"                   # D4 (new shared string)
$ws.Cells.Item(4,5).Value  = 3                      # E4
$ws.Cells.Item(4,6).Value  = 1                      # F4
$ws.Cells.Item(4,7).Value  = 1.995314333333334      # G4
$ws.Cells.Item(4,8).Value  = 5.985943000000001      # H4
$ws.Cells.Item(4,9).Value  = 1                      # I4
$ws.Cells.Item(4,10).Value = 1                      # J4
$ws.Cells.Item(4,11).Value = 1                      # K4
$ws.Cells.Item(4,12).Value = 0.3333333333333333     # L4
$ws.Cells.Item(4,13).Value = 0.004340666666666667   # M4
$ws.Cells.Item(4,14).Value = 0.013022               # N4
$ws.Cells.Item(4,15).Value = 0.0005224795816889323  # O4
$ws.Cells.Item(4,16).Value = 0.0005224795816889323  # P4
$ws.Cells.Item(4,17).Value = 0.008660994416222225   # Q4
$ws.Cells.Item(4,18).Value = 0.07794894974600002    # R4
$ws.Cells.Item(4,19).Value = 0.0005224795816889323  # S4
$ws.Cells.Item(4,20).Value = 0.0005224795816889323  # T4

# --- New row 5 (sCs / Gdnf / Ret / sCs), containing what used to be row 4's numbers (updated) ---
$ws.Cells.Item(5,1).Value  = "sCs"                  # A5
$ws.Cells.Item(5,2).Value  = "Gdnf"                 # B5
$ws.Cells.Item(5,3).Value  = "Ret"                  # C5
$ws.Cells.Item(5,4).Value  = "sCs"                  # D5
$ws.Cells.Item(5,5).Value  = 3                      # E5
$ws.Cells.Item(5,6).Value  = 1                      # F5
$ws.Cells.Item(5,7).Value  = 1.995314333333334      # G5
$ws.Cells.Item(5,8).Value  = 5.985943000000001      # H5
$ws.Cells.Item(5,9).Value  = 1                      # I5
$ws.Cells.Item(5,10).Value = 1                      # J5
$ws.Cells.Item(5,11).Value = 3                      # K5
$ws.Cells.Item(5,12).Value = 1                      # L5
$ws.Cells.Item(5,13).Value = 0.7855989999999999     # M5
$ws.Cells.Item(5,14).Value = 2.356797               # N5
$ws.Cells.Item(5,15).Value = 0.09456138156087625    # O5
$ws.Cells.Item(5,16).Value = 0.09456138156087625    # P5
$ws.Cells.Item(5,17).Value = 1.567516944952333      # Q5
$ws.Cells.Item(5,18).Value = 14.107652504571        # R5
$ws.Cells.Item(5,19).Value = 0.09456138156087625    # S5
$ws.Cells.Item(5,20).Value = 0.09456138156087625    # T5
